$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Column D price cells to be stored as Text so that numeric-looking
# strings (e.g. "1.00", "0.998") keep their exact literal representation
# instead of being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.655.18"
$ws.Range("E2").Value = "  -7.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.686.13"
$ws.Range("E3").Value = "  -6.93%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.05"
$ws.Range("E5").Value = "  -6.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.70"
$ws.Range("E6").Value = "  +4.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.677.00"
$ws.Range("E7").Value = "  -6.96%  "
$ws.Range("E8").Value = "  -7.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.712"
$ws.Range("E10").Value = "  -5.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.164"
$ws.Range("E11").Value = "  -10.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.05"
$ws.Range("E12").Value = "  -7.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000298"
$ws.Range("E13").Value = "  -10.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.65"
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.288.05"
$ws.Range("E15").Value = "  -6.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.725.25"
$ws.Range("E16").Value = "  -6.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.36"
$ws.Range("E17").Value = "  -5.39%  "
$ws.Range("E18").Value = "  -3.24%  "
$ws.Range("B19").Value = "Polygon"
$ws.Range("C19").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.14"
$ws.Range("E19").Value = "  -9.07%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.94"
$ws.Range("E20").Value = "  -8.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.768.60"
$ws.Range("E21").Value = "  -6.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "408.35"
$ws.Range("E22").Value = "  -6.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.54"
$ws.Range("E23").Value = "  -7.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.08"
$ws.Range("E24").Value = "  -8.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.08"
$ws.Range("E25").Value = "  -8.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.81"
$ws.Range("E26").Value = "  -9.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.71"
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.83"
$ws.Range("E28").Value = "  -5.89%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.49"
$ws.Range("E30").Value = "  -9.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.93"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.93"
$ws.Range("E32").Value = "  -8.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.69"
$ws.Range("E33").Value = "  -6.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.117"
$ws.Range("E34").Value = "  -9.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.20"
$ws.Range("E35").Value = "  -7.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.63"
$ws.Range("E36").Value = "  -6.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0914"
$ws.Range("E37").Value = "  -12.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "601.15"
$ws.Range("E38").Value = "  -7.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.400"
$ws.Range("E39").Value = "  -7.17%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.26"
$ws.Range("E41").Value = "  +13.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  -7.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.06"
$ws.Range("E44").Value = "  -11.51%  "
$ws.Range("E45").Value = "  -9.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.61"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.45"
$ws.Range("E47").Value = "  -11.33%  "
$ws.Range("E48").Value = "  -9.31%  "
$ws.Range("E49").Value = "  -15.58%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000268"
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.733.72"
$ws.Range("E51").Value = "  -3.53%  "
